# "update with single HB session to lbs side to test current state"
#
# Appends one new climbing-session row (row 57) to Sheet1, mirroring the
# existing table of sessions (date, jug, imr, med_edge, pinch_med, mr_2fp,
# large_edge, pinch_wide, sloper, mrp_3fp), widens column E a touch to fit
# the new (longer) value, and leaves the view scrolled down with the newly
# added cell selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new session row (row 57) ---------------------------------------------
$ws.Range("A57").Value = "9 Nov 2025"
$ws.Range("B57").Value = "0,6"
$ws.Range("C57").Value = "-7.5,5,9"
$ws.Range("D57").Value = "-17.5,3,5,7,5"
$ws.Range("E57").Value = "-50,0,4,3,3,2,2,1"
$ws.Range("F57").Value = "-30,2,3,5,3,3"
$ws.Range("G57").Value = "-15,2,9,4,4,5"
$ws.Range("H57").Value = "-50,2,6,5,3,2"
$ws.Range("I57").Value = "-45,4,7,6"
$ws.Range("J57").Value = "-35,2,7,6,4,4"

# --- column E grew a bit wider to fit "-50,0,4,3,3,2,2,1" ------------------
$ws.Columns.Item(5).ColumnWidth = 13.66

# --- scroll the view down and leave the new row's E cell selected ----------
$ws.Range("E57").Select()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
